# edit.ps1 - applies the Deckblatt.docx changes described by the diff:
#   1. "MAD – Network Monitoring" -> "LAN – " + "Monitoring" (two runs, same formatting)
#   2. "Höhere ... Bundeslehr- und Versuchsanstalt ..." split into two runs with a
#      "_GoBack" bookmark inserted between them
#   3. remove the old "_GoBack" bookmark that used to sit in the "Ausgeführt im
#      Schuljahr 2014/15 von:" table cell (Word always keeps only one _GoBack,
#      relocated to the most recent edit point)

$d = $word.ActiveDocument
$enDash = [char]0x2013
$szlig  = [char]0x00DF

# -----------------------------------------------------------------
# Change 1: "MAD <enDash> Network Monitoring" -> "LAN <enDash> Monitoring"
#           split into "LAN <enDash> " / "Monitoring" runs
# -----------------------------------------------------------------

# 1a. Text substitution (stays a single run, formatting untouched)
$r1 = $d.Content
$null = $r1.Find.Execute("MAD", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "LAN", 2)

$r1b = $d.Content
$null = $r1b.Find.Execute("Network ", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 2)

# 1b. Force a run split right before "Monitoring" (re-apply identical Bold
#     value so formatting is unchanged but the run boundary is created)
$r1c = $d.Content
$null = $r1c.Find.Execute("Monitoring", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$bold1 = $r1c.Font.Bold
$r1c.Font.Bold = -not $bold1
$r1c.Font.Bold = $bold1

# -----------------------------------------------------------------
# Change 2: split "Höhere ... Bundeslehr- und Versuchsanstalt ..." into two
#           runs and drop a "_GoBack" bookmark at the split point
# -----------------------------------------------------------------

$r2 = $d.Content
$target2 = "und Versuchsanstalt Anichstra" + $szlig + "e"
$null = $r2.Find.Execute($target2, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$bold2 = $r2.Font.Bold
$r2.Font.Bold = -not $bold2
$r2.Font.Bold = $bold2

$bmRange = $d.Range($r2.Start, $r2.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -----------------------------------------------------------------
# Change 3: remove the previous "_GoBack" bookmark that sat between the
#           "5" run and the " von:   " run, without disturbing those runs'
#           text/formatting. A Find/Replace spanning the bookmark clears it,
#           then we re-split the merged run back into its original pieces.
# -----------------------------------------------------------------

$r3 = $d.Content
$null = $r3.Find.Execute("5 ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "5 ", 2)

# restore "5" as its own run (split off the trailing " von:   ")
$r3b = $d.Content
$null = $r3b.Find.Execute("5 von:   ", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$bold3 = $r3b.Font.Bold
$r3b.Font.Bold = -not $bold3
$r3b.Font.Bold = $bold3

# restore " von:   " as its own run (split off from "5")
$r3c = $d.Content
$null = $r3c.Find.Execute(" von:   ", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$bold3c = $r3c.Font.Bold
$r3c.Font.Bold = -not $bold3c
$r3c.Font.Bold = $bold3c

Write-Output "done"
